# feat: add 2022-Q3 data
#
# 1. Duplicate the "2022-Q1" sheet (this clones formatting/styles/dimension
#    exactly) as a template for the new "2022-Q3" sheet, inserted right
#    before "2022-Q1" so the tab order becomes:
#      总计, 2022-Q3, 2022-Q1, 2021-Q2
$wb = $excel.ActiveWorkbook

$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1)
$q3 = $wb.Worksheets.Item("2022-Q1 (2)")
$q3.Name = "2022-Q3"

# --- Populate the new "2022-Q3" sheet -------------------------------------
# Row1 (headers) and row2 A/B/H already match "2022-Q1" after the clone, so
# only the differing cells need to be written.
$q3.Range("C2").Value = "上投摩根亚太优势混合（QDII）"
$q3.Range("H2").Value = 8

# Numeric-looking text cells: write through a staging cell so the engine
# stores them as literal text (matching the source workbook's convention of
# keeping these metric columns as text) instead of auto-coercing to a
# number, then copy only the value (no formatting) into place.
function Set-TextValue($ws, $cellAddr, $text) {
    $ws.Range("Z100").NumberFormat = "@"
    $ws.Range("Z100").Value = $text
    $ws.Range("Z100").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $ws.Range("Z100").Clear()
}

Set-TextValue $q3 "D2" "22.76"
Set-TextValue $q3 "E2" "90.69"
Set-TextValue $q3 "F2" "2.30"
Set-TextValue $q3 "G2" "0.5235"

# Row3 is brand new (the "2022-Q1" template only had one data row) - build
# it by copying row2's formatting down (so A3 picks up style "s=2" the same
# way A2 has it) and then overwriting the values.
$q3.Range("A2:H2").Copy()
$q3.Range("A3").PasteSpecial(-4122)

$q3.Range("A3").Value = 1
Set-TextValue $q3 "B3" "006105"
$q3.Range("C3").Value = "泰达宏利印度机会股票（QDII）"
Set-TextValue $q3 "D3" "0.66"
Set-TextValue $q3 "E3" "86.08"
Set-TextValue $q3 "F3" "2.60"
Set-TextValue $q3 "G3" "0.0172"
$q3.Range("H3").Value = 8

# --- Update the "总计" summary sheet ---------------------------------------
# Insert the 2022-Q3 row after the header and push the existing two rows
# down by one. Values are written explicitly (rather than via Rows.Insert)
# to avoid inheriting stray formatting/number noise from a shift operation.
$zj = $wb.Worksheets.Item("总计")

# Row4 <= old Row3 (2021-Q2), written first so nothing is clobbered while
# row2/row3 are still being read below.
$zj.Range("A2:D2").Copy()
$zj.Range("A4").PasteSpecial(-4122)
$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2021-Q2"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 0.78

# Row3 <= old Row2 (2022-Q1)
$zj.Range("B3").Value = "2022-Q1"
$zj.Range("C3").Value = 1
$zj.Range("D3").Value = 0.5

# Row2 <= new 2022-Q3 summary
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0.54

# Keep "2021-Q2" as the selected tab, matching the source workbook.
$wb.Worksheets.Item("2021-Q2").Activate()
